$wb = $excel.ActiveWorkbook

# --- Sheet "NPCs": add new row for LabyrinthWeaver ---
$npcs = $wb.Worksheets.Item("NPCs")

$npcs.Cells.Item(10, 1).Value = "LabyrinthWeaver"
$npcs.Cells.Item(10, 2).Value = "Labyrinth Weaver"
$npcs.Cells.Item(10, 3).Value = 2
$npcs.Cells.Item(10, 4).Value = "Labyrinth"
$npcs.Cells.Item(10, 6).Value = 1
$npcs.Cells.Item(10, 7).Value = "/m LabyrinthWeaver:"
$npcs.Cells.Item(10, 8).Value = 384
$npcs.Cells.Item(10, 9).Value = 496

# --- Sheet "Npcs Commands": add new row for the Labyrinth Weaver's command ---
$cmds = $wb.Worksheets.Item("Npcs Commands")

$cmds.Cells.Item(10, 1).Value = "Labyrinth Weaver"
$cmds.Cells.Item(10, 2).Value = "Make Sash"
$cmds.Cells.Item(10, 3).Value = 0
